$d = $word.ActiveDocument

# Task E ("9-11 task E") input-data limit: change "n <= 10**5)." to
# "n <= 10**6)." — i.e. only the exponent digit (5 -> 6) changes.

# 1) Plain text substitution for the whole "10**5)." -> "10**6)." span.
$rng = $d.Content
$rng.Find.Execute("10**5).", $true, $false, $false, $false, $false, $true, 1, $false, "10**6).", 2)

# 2) Re-find the (now updated) text so we can isolate just the changed
#    digit "6" and force Word to split it into its own run, matching
#    how a real edit (select the "5", type "6") would leave the OOXML:
#    the surrounding text stays in neighboring runs while the replaced
#    character gets its own run.
$rng2 = $d.Content
$rng2.Find.Execute("10**6).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start
# "10**6)." -> indices: 1,0,*,*,6,),.  -> the digit "6" is at offset 4
$digitRange = $d.Range($start2 + 4, $start2 + 5)
$digitRange.Font.Bold = $true
$digitRange.Font.Bold = $false
